$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a pure number but must stay text
# (matches the source data which stores everything as inline strings).
# Temporarily mark them as text before assigning, then restore the default
# "Normal" style so no stray number-format style is left on the cell.
$textForceCells = @(
    "D5",
    "D6",
    "D7",
    "D10",
    "D11",
    "D12",
    "D16",
    "D19",
    "D20",
    "D21",
    "D22",
    "D24",
    "D25",
    "D27",
    "D28",
    "D29",
    "D30",
    "D32",
    "D33",
    "D35",
    "D36",
    "D37",
    "D42",
    "D43",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51",
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '64.422.47'
$ws.Range("E2").Value = '  -2.52%  '
$ws.Range("D3").Value = '3.176.32'
$ws.Range("E3").Value = '  -4.16%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '570.95'
$ws.Range("E5").Value = '  -2.57%  '
$ws.Range("D6").Value = '168.67'
$ws.Range("E6").Value = '  -7.49%  '
$ws.Range("D7").Value = '0.607'
$ws.Range("E7").Value = '  -7.14%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").Value = '3.185.61'
$ws.Range("E9").Value = '  -3.84%  '
$ws.Range("D10").Value = '0.120'
$ws.Range("E10").Value = '  -4.27%  '
$ws.Range("D11").Value = '6.81'
$ws.Range("E11").Value = '  -0.19%  '
$ws.Range("D12").Value = '0.388'
$ws.Range("E12").Value = '  -3.24%  '
$ws.Range("D13").Value = '3.736.47'
$ws.Range("E13").Value = '  -3.90%  '
$ws.Range("E14").Value = '  -1.38%  '
$ws.Range("D15").Value = '64.491.55'
$ws.Range("E15").Value = '  -2.49%  '
$ws.Range("D16").Value = '25.27'
$ws.Range("E16").Value = '  -3.48%  '
$ws.Range("E17").Value = '  -3.90%  '
$ws.Range("D18").Value = '3.184.24'
$ws.Range("E18").Value = '  -4.05%  '
$ws.Range("D19").Value = '418.51'
$ws.Range("E19").Value = '  -1.49%  '
$ws.Range("D20").Value = '12.98'
$ws.Range("E20").Value = '  -1.21%  '
$ws.Range("D21").Value = '5.36'
$ws.Range("E21").Value = '  -3.28%  '
$ws.Range("D22").Value = '7.12'
$ws.Range("E22").Value = '  -3.50%  '
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("D24").Value = '70.21'
$ws.Range("E24").Value = '  -2.02%  '
$ws.Range("D25").Value = '5.67'
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("E26").Value = '  +2.20%  '
$ws.Range("D27").Value = '0.495'
$ws.Range("E27").Value = '  -3.26%  '
$ws.Range("D28").Value = '0.0000105'
$ws.Range("E28").Value = '  -7.21%  '
$ws.Range("D29").Value = '8.71'
$ws.Range("E29").Value = '  -1.82%  '
$ws.Range("D30").Value = '0.997'
$ws.Range("E30").Value = '  -0.41%  '
$ws.Range("E31").Value = '  -3.69%  '
$ws.Range("D32").Value = '21.76'
$ws.Range("E32").Value = '  -2.75%  '
$ws.Range("D33").Value = '0.998'
$ws.Range("E33").Value = '  -0.12%  '
$ws.Range("E34").Value = '  -2.26%  '
$ws.Range("D35").Value = '6.35'
$ws.Range("E35").Value = '  -3.17%  '
$ws.Range("D36").Value = '1.12'
$ws.Range("E36").Value = '  -4.07%  '
$ws.Range("D37").Value = '156.83'
$ws.Range("E37").Value = '  -2.24%  '
$ws.Range("E38").Value = '  -5.08%  '
$ws.Range("E39").Value = '  -5.26%  '
$ws.Range("D40").Value = '2.689.54'
$ws.Range("E40").Value = '  -6.39%  '
$ws.Range("E41").Value = '  -1.90%  '
$ws.Range("D42").Value = '24.25'
$ws.Range("E42").Value = '  -7.94%  '
$ws.Range("D43").Value = '39.31'
$ws.Range("E43").Value = '  -1.31%  '
$ws.Range("E44").Value = '  -5.57%  '
$ws.Range("D45").Value = '0.0623'
$ws.Range("E45").Value = '  -5.74%  '
$ws.Range("D46").Value = '5.58'
$ws.Range("E46").Value = '  -5.38%  '
$ws.Range("D47").Value = '0.0263'
$ws.Range("E47").Value = '  -2.96%  '
$ws.Range("D48").Value = '291.44'
$ws.Range("E48").Value = '  -7.14%  '
$ws.Range("D49").Value = '21.33'
$ws.Range("E49").Value = '  -7.54%  '
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").Value = '0.0993'
$ws.Range("E50").Value = '  -6.29%  '
$ws.Range("B51").Value = 'FirstDigitalUSD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D51").Value = '0.998'
$ws.Range("E51").Value = '  -0.15%  '

foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}

